$d = $word.ActiveDocument

# Locate the target paragraph: "Salidas= a.b'.b.c + a.b.c"
$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "Salidas=*") {
        $target = $para
    }
}

$r = $target.Range

# Find the last occurrence of "a.b.c" in that paragraph (the un-bookmarked one at the end)
$found = $r.Find.Execute("a.b.c", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Duplicate the found run's formatting (keeps sz/szCs that plain property
# assignment can't reach) by copying its FormattedText into a fresh spot
# right after it, then fixing up the copied text in place.
$srcRange = $r.Duplicate
$ft = $srcRange.FormattedText

$insertPoint = $r.End
$r.Collapse(0)
$r.FormattedText = $ft

$newRange = $d.Range($insertPoint, $insertPoint + 5)
$newRange.Find.Execute("a.b.c", $false, $false, $false, $false, $false, $true, 1, $false, " + á.b.c", 2) | Out-Null

# Move the _GoBack bookmark from its old spot (mid "b.´b.c") to the very
# end of the paragraph, right after the text we just appended.
$paraEnd = $target.Range.End - 1
$bmRange = $d.Range($paraEnd, $paraEnd)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
